$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "Jasmine Scottini"
$ws.Range("B37").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C37").Value = "Federico Fasanelli | SBARX"
$ws.Range("D37").Value = "Filippo Benetti | I Magnifici"
$ws.Range("E37").Value = "Luca Frasca | Clitoriders"
$ws.Range("F37").Value = "Francesco Cristoforetti | Vigili del Fusto"
